$d = $word.ActiveDocument

# Split paragraph 1 (which currently holds "1" plus the trailing
# _GoBack bookmark) right after the "1" run, i.e. right before the
# bookmark, so the bookmark carries over into the newly created
# second paragraph.
$splitPoint = $d.Range(1, 1)
$splitPoint.InsertParagraphAfter()

# Insert the new paragraph's text. Append a throwaway trailing
# character along with "333" so the insertion point is a real
# (non-collapsed) range - this runtime mis-handles Bookmarks.Add on a
# range that is collapsed exactly at a paragraph-end position. We'll
# trim the extra character back off afterwards, which correctly
# collapses the bookmark to that now-empty spot.
$para2 = $d.Paragraphs(2)
$insertRange = $para2.Range
$insertRange.InsertBefore("333X")

$para2b = $d.Paragraphs(2)
$r2b = $para2b.Range
$lastChar = $d.Range($r2b.Start + 3, $r2b.Start + 4)

# Re-anchor the _GoBack bookmark onto the placeholder character.
$d.Bookmarks.Add("_GoBack", $lastChar)

# Delete the placeholder character; the bookmark collapses onto the
# now-empty range right after "333" and before the paragraph mark.
$lastChar2 = $d.Range($r2b.Start + 3, $r2b.Start + 4)
$lastChar2.Text = ""
